$d = $word.ActiveDocument

# 1 & 2. Title (in header table) and the "Congratulations" paragraph both contain the
#        same old sentence fragment; wdReplaceAll (2) over the whole document content
#        updates both occurrences in one pass.
$null = $d.Content.Find.Execute(
    "Ask the court to issue a protection order or order granting you possession relative to domesticated animals",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "File a Petition 209 Relative to Domestic Animals", 2)

# 3. "Tell the judge ..." paragraph: statute citation replaced with the new phrase.
$null = $d.Content.Find.Execute(
    "Petition Filed Pursuant to G. L. c. 209A, Sect. 11 Relative To Domesticated Animal(s)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Petition209as11relativetodomesticatedanimalssecond 2 2", 2)

# 4. Remove the stray "_GoBack" bookmark (and its now-empty wrapper paragraph
#    collapses to a bare <w:p/>).
$gb = $d.Bookmarks.Item("_GoBack")
$gb.Delete()
